$d = $word.ActiveDocument

# Helper: replace the contents of a paragraph (everything except its
# trailing paragraph mark, so the pPr / numbering / paragraph-level rPr
# stay untouched) with a literal run of WordprocessingML.
function Set-ParaInnerXml($paraIndex, [string]$innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    $inner = $d.Range($full.Start, $full.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $inner.InsertXML($pkg)
}

# --- Change 1: "Make Stacked Bar Graphs" paragraph gets a left tab stop at 5130 twips (256.5 pt) ---
$p5 = $d.Paragraphs.Item(5)
$p5.Range.ParagraphFormat.TabStops.Add(256.5)

# --- Change 2: "To-DO:" paragraph - split "To-DO" off into its own
#     spell-checked run, leaving ":" as a separate trailing run ---
$xml14 = '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>To-DO</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>' +
         '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>:</w:t></w:r>'
Set-ParaInnerXml 14 $xml14

# --- Change 3: "DONE – Solar_All_County" paragraph - "DONE – " becomes
#     its own run and "Solar_All_County" (3 runs) is wrapped in proofErr ---
$xml22 = '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">DONE – </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Solar_</w:t></w:r>' +
         '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>All</w:t></w:r>' +
         '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>_County</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml 22 $xml22

# --- Change 4: "DONE – Solar_Rates_County" paragraph - same split pattern ---
$xml30 = '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t xml:space="preserve">DONE – </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:color w:val="00B050"/></w:rPr><w:t>Solar_Rates_County</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml 30 $xml30

# --- Change 5: "Table is called Solar_Rates_County" paragraph - keep
#     <w:lastRenderedPageBreak/> on the first run, split the table name off ---
$xml34 = '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Table is called </w:t></w:r>' +
         '<w:proofErr w:type="spellStart"/>' +
         '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>Solar_Rates_County</w:t></w:r>' +
         '<w:proofErr w:type="spellEnd"/>'
Set-ParaInnerXml 34 $xml34

Write-Host "All edits applied"
